$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 - Getriebemotor 12 V (Conrad)
$ws.Range("C12").Value = "Getriebemotor 12 V"
$ws.Range("D12").Value = "IG320005-3AC21R"
$ws.Range("E12").Value = "Conrad"
$ws.Range("F12").Value = "234253 - UP"
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 27.95

# Row 13 - 4S LiPo Akku 500mAh (Hobbyking)
$ws.Range("C13").Value = "4S LiPo Akku 500mAh"
$ws.Range("E13").Value = "Hobbyking"
$ws.Range("F13").Value = "9067000341-0"
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 30.12

# Row 14 - 4mm Gold-Steckverbinder 10 Stk. (Hobbyking)
$ws.Range("C14").Value = "4mm Gold-Steckverbinder 10 Stk."
$ws.Range("E14").Value = "Hobbyking"
$ws.Range("F14").Value = 15000068
$ws.Range("F14").HorizontalAlignment = -4131
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 3.16

# Row 15 - Lipoly Niederspannungs-Warnung (2s ~ 4s) (Hobbyking)
$ws.Range("C15").Value = "Lipoly Niederspannungs-Warnung (2s ~ 4s)"
$ws.Range("E15").Value = "Hobbyking"
$ws.Range("F15").Value = "DL-Volt-Alarm"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 2.63

# Row 16 - Schleifring (servo technica)
$ws.Range("C16").Value = "Schleifring"
$ws.Range("D16").Value = "SVTS C 03-X-A-00/06"
$ws.Range("E16").Value = "servo technica"
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 1

# Widen column C to fit the longer descriptions now stored there
$ws.Columns("C").ColumnWidth = 35.59

# Move the active selection, matching the saved view state
$ws.Range("E20").Select()
